$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the current selection/active cell to E4
$ws.Range("E4").Select()

# Increase the custom height of row 3 to 45
$ws.Rows.Item(3).RowHeight = 45
